# Add season record columns (Wins / Losses / Ties) to the player table.
# The old scraper only grabbed team statistics, not the season record, so
# three new columns are appended after the existing data (which ends at AC)
# and populated with the team's record for every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, formatted like the rest of the header row (AD1:AF1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Season record values for every player row (rows 2 through 46).
$wins = 91
$losses = 70
$ties = 0

for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
